# Fix mistake in IndBeskr data: IndID and Register were swapped for nra.
# Swap the values of column A (IndID) and column B (Register) for every
# row belonging to the "nra" register (rows 39-51 on the
# "Indikatorbeskrivelser" sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indikatorbeskrivelser")

for ($r = 39; $r -le 51; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    $aVal = $aCell.Value2
    $bVal = $bCell.Value2

    $aCell.Value2 = $bVal
    $bCell.Value2 = $aVal
}
